# Generate Report for Handoff
# File "7332775b-08c9-4371-bbb2-94e10f68287a.md" has moved from "In Translation"
# to "Ready for handoff" - update the status + timestamp columns on all three
# sheets (Overview, zh-cn, de-de) for that file's row (row 7).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E7").Value = "Ready for handoff"
$overview.Range("F7").Value = "Ready for handoff"
$overview.Range("G7").Value = "2016-08-21 14:52:06"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C7").Value = "Ready for handoff"
$zhcn.Range("H7").Value = "2016-08-21 14:51:58"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C7").Value = "Ready for handoff"
$dede.Range("H7").Value = "2016-08-21 14:52:06"
